{"js": "// Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact\n// metrics (percentages, dollar amounts, large numbers) across the resume.\n//\n// Strategy: for each paragraph that contains metric(s) to highlight, run a\n// paragraph-scoped search for the exact metric substring (search is scoped\n// to the paragraph so repeated metric strings elsewhere in the document are\n// left untouched) and apply bold + the corporate-blue color to every\n// matching range. Metrics are processed in left-to-right order per\n// paragraph so that, e.g., \"23%\" is bolded before \"64%\" in the same\n// sentence \u2014 this mirrors the run-splitting seen in the target OOXML.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of paragraph text (exact, pre-edit) -> ordered list of metric\n// substrings inside that paragraph that must become bold + colored.\nconst metricPlan = [\n  {\n    text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"]\n  },\n  {\n    text: \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    metrics: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    metrics: [\"$2\"]\n  },\n  {\n    text: \"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n    metrics: [\"57%\"]\n  },\n  {\n    text: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    metrics: [\"73.5%\"]\n  },\n  {\n    text: \"\u2022 $4.7M savings enabled nonprofit access\",\n    metrics: [\"$4.7M\"]\n  },\n  {\n    text: \"\u2022 178% accuracy improvement in racial classification algorithms\",\n    metrics: [\"178%\"]\n  }\n];\n\n// Build a quick lookup from trimmed paragraph text to its metric list so we\n// don't depend on a fixed paragraph index (robust to minor doc changes).\nfunction normalize(t) {\n  return (t || \"\").replace(/\\r/g, \"\").trim();\n}\n\nconst planByText = new Map();\nfor (const entry of metricPlan) {\n  planByText.set(normalize(entry.text), entry.metrics);\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const key = normalize(para.text);\n  const metrics = planByText.get(key);\n  if (!metrics) continue;\n\n  for (const metric of metrics) {\n    const results = para.search(metric, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const range of results.items) {\n      range.font.bold = true;\n      range.font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact\n# metrics (percentages, dollar amounts, large numbers) across the resume.\n#\n# Strategy: for each paragraph whose full text matches one of the target\n# sentences, run Find.Execute scoped to a fresh Range built from that\n# paragraph for each metric substring (in left-to-right order) and apply\n# Bold + the corporate-blue font color to the matched text. Re-fetching\n# $paragraph.Range for every search keeps the Find call properly scoped to\n# the paragraph even after earlier matches have split it into more runs.\n\nfunction Get-WdColor([int]$r, [int]$g, [int]$b) {\n    # Word's Font.Color is a 0x00BBGGRR integer (not plain RGB).\n    return ($b * 65536) + ($g * 256) + $r\n}\n\n$highlightColor = Get-WdColor 44 62 80   # #2C3E50\n\n# Ordered metric substrings to bold+color, keyed by the paragraph's exact\n# (trimmed) text before editing.\n$plan = @(\n    @{\n        Text = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text = \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\"\n        Metrics = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\")\n    },\n    @{\n        Text = \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        Metrics = @(\"73.5%\", '$4.7M')\n    },\n    @{\n        Text = \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\"\n        Metrics = @('$2')\n    },\n    @{\n        Text = \"\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\"\n        Metrics = @(\"57%\")\n    },\n    @{\n        Text = \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n        Metrics = @(\"73.5%\")\n    },\n    @{\n        Text = \"\u2022 `$4.7M savings enabled nonprofit access\"\n        Metrics = @('$4.7M')\n    },\n    @{\n        Text = \"\u2022 178% accuracy improvement in racial classification algorithms\"\n        Metrics = @(\"178%\")\n    }\n)\n\n# Build a lookup keyed by paragraph text for quick matching.\n$planByText = @{}\nforeach ($entry in $plan) {\n    $planByText[$entry.Text] = $entry.Metrics\n}\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $paraText = $para.Range.Text.TrimEnd([char]13, [char]10)\n\n    if (-not $planByText.ContainsKey($paraText)) {\n        continue\n    }\n\n    $metrics = $planByText[$paraText]\n    foreach ($metric in $metrics) {\n        $rng = $para.Range\n        $rng.Find.MatchCase = $true\n        $rng.Find.MatchWildcards = $false\n        $found = $rng.Find.Execute($metric)\n        if ($found) {\n            $rng.Font.Bold = 1\n            $rng.Font.Color = $highlightColor\n        }\n    }\n}\n"}
